$d = $word.ActiveDocument

# --- Change 3: remove the blank "Historia n" template block (5 paragraphs) ---
$rngFind = $d.Content
$found = $rngFind.Find.Execute("Historia n: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "anchor Historia n: not found" }
$pHistoriaN = $rngFind.Paragraphs(1)
$pBlank = $pHistoriaN.Previous()
$pComo = $pHistoriaN.Next()
$pQuiero = $pComo.Next()
$pPara = $pQuiero.Next()
$delRange = $d.Range($pBlank.Range.Start, $pPara.Range.End)
$delRange.Delete()

# --- Change 2: Historia 7 Como/Quiero/Para - replace shd highlight with plain white highlight ---
$rngFind2 = $d.Content
$found2 = $rngFind2.Find.Execute("Quiero: ver los trabajos e investigaciones que ha desarrollado el docente.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "anchor Quiero: ver los trabajos not found" }
$pQuieroH7 = $rngFind2.Paragraphs(1)
$pComoH7 = $pQuieroH7.Previous()
$pParaH7 = $pQuieroH7.Next()
$blockRange = $d.Range($pComoH7.Range.Start, $pParaH7.Range.End)
$xmlFrag2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:ind w:left="1440" w:hanging="360"/><w:rPr><w:highlight w:val="white"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="white"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Como: estudiante</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:ind w:left="1440" w:hanging="360"/><w:rPr><w:highlight w:val="white"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="white"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Quiero: ver los trabajos e investigaciones que ha desarrollado el docente.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:ind w:left="1440" w:hanging="360"/><w:rPr><w:highlight w:val="white"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="white"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Para: escoger un docente que tenga un campo de desarrollo afín con los gustos del estudiante a la hora de matricular una asignatura o ingresar a uno de los grupos de investigación en el que el docente participa de alguna manera, además de la posibilidad de saber quién ha trabajado en temas afines a la tesis elegida .</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$blockRange.InsertXML($xmlFrag2)

# --- Change 1: merge "-Profesores" paragraph with "-Directivos.(?)" paragraph, clearing the latter text ---
$rngFind3 = $d.Content
$found3 = $rngFind3.Find.Execute("-Profesores", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "anchor -Profesores not found" }
$pProfesores = $rngFind3.Paragraphs(1)
$pDirectivos = $pProfesores.Next()
$mergeRange = $d.Range($pProfesores.Range.Start, $pDirectivos.Range.End)
$xmlFrag1 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">-Profesores</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$mergeRange.InsertXML($xmlFrag1)

Write-Output "Done"
